$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all wt_level values (column A, rows 2 to 74) to 0, clarifying that wt has only 1 level.
$ws.Range("A2:A74").Value = 0

# Update the view's selection / scroll position to match the final saved state.
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("F8").Select()
